$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 363
$ws.Range("G2").Value = 3708
$ws.Range("G8").Value = 7522
$ws.Range("G9").Value = 8416
$ws.Range("I9").Value = 17076
$ws.Range("G10").Value = 9142
$ws.Range("G11").Value = 10501
$ws.Range("I11").Value = 20839
$ws.Range("G12").Value = 12140
$ws.Range("G13").Value = 14072
$ws.Range("I13").Value = 26445
$ws.Range("G14").Value = 14488
$ws.Range("G15").Value = 15760
$ws.Range("G16").Value = 16450
$ws.Range("G18").Value = 17870
$ws.Range("I18").Value = 34902
$ws.Range("G19").Value = 17230
$ws.Range("C20").Value = 1061
$ws.Range("E21").Value = 5332
$ws.Range("G21").Value = 23818
$ws.Range("G22").Value = 27016
$ws.Range("G23").Value = 30552
$ws.Range("E24").Value = 6420
$ws.Range("G24").Value = 33350
$ws.Range("I24").Value = 53867
$ws.Range("G25").Value = 33461
$ws.Range("G26").Value = 35886
$ws.Range("G27").Value = 41805
$ws.Range("I27").Value = 69619
$ws.Range("G28").Value = 46335
$ws.Range("I28").Value = 75164
$ws.Range("G29").Value = 49929
$ws.Range("I29").Value = 81366
$ws.Range("G30").Value = 52222
$ws.Range("I30").Value = 86556
$ws.Range("G31").Value = 55529
$ws.Range("I31").Value = 93024
$ws.Range("G32").Value = 57810
$ws.Range("I32").Value = 97805
$ws.Range("G33").Value = 59274
$ws.Range("I33").Value = 100420
$ws.Range("C34").Value = 3031
$ws.Range("G34").Value = 64036
$ws.Range("I34").Value = 108835
$ws.Range("C35").Value = 3320
$ws.Range("E35").Value = 19026
$ws.Range("G35").Value = 68606
$ws.Range("I35").Value = 116852
$ws.Range("C36").Value = 5098
$ws.Range("E36").Value = 28387
$ws.Range("G36").Value = 70192
$ws.Range("I36").Value = 118726
